$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Activate the "Repayment schedule" sheet (it becomes the active/selected tab)
$ws.Activate()

# Insert a new blank column before column N (14th column), shifting the
# existing N..P columns (Late / heading / Outstanding) one to the right.
$ws.Columns("N:N").Insert()

# The newly inserted column inherits its width from the column to its left
# (column M, "In Advance").
$ws.Columns("N:N").ColumnWidth = $ws.Columns("M:M").ColumnWidth

# Leave the final selection on the sheet where the edit was made.
$ws.Range("R6").Select()
